$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "68.273.84"
$ws.Range("E2").Value = "  -2.99%  "
$ws.Range("D3").Value = "3.390.46"
$ws.Range("E3").Value = "  -3.67%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "607.03"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "164.43"
$ws.Range("E6").Value = "  -5.99%  "
$ws.Range("D7").Value = "0.588"
$ws.Range("E7").Value = "  -3.72%  "
$ws.Range("D8").Value = "3.384.35"
$ws.Range("E8").Value = "  -3.71%  "
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").Value = "0.189"
$ws.Range("E10").Value = "  -2.65%  "
$ws.Range("D11").Value = "6.79"
$ws.Range("E11").Value = "  -6.17%  "
$ws.Range("D12").Value = "0.548"
$ws.Range("E12").Value = "  -5.96%  "
$ws.Range("D13").Value = "43.11"
$ws.Range("E13").Value = "  -6.79%  "
$ws.Range("D14").Value = "0.0000263"
$ws.Range("E14").Value = "  -4.52%  "
$ws.Range("D15").Value = "3.963.41"
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("D16").Value = "7.96"
$ws.Range("E16").Value = "  -4.32%  "
$ws.Range("D17").Value = "3.418.10"
$ws.Range("E17").Value = "  -2.70%  "
$ws.Range("D18").Value = "68.514.64"
$ws.Range("E18").Value = "  -2.64%  "
$ws.Range("D19").Value = "568.66"
$ws.Range("E19").Value = "  -7.31%  "
$ws.Range("D20").Value = "0.119"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").Value = "16.77"
$ws.Range("E21").Value = "  -2.96%  "
$ws.Range("D22").Value = "0.829"
$ws.Range("E22").Value = "  -5.32%  "
$ws.Range("D23").Value = "8.78"
$ws.Range("E23").Value = "  -3.15%  "
$ws.Range("D24").Value = "93.56"
$ws.Range("E24").Value = "  -3.89%  "
$ws.Range("D25").Value = "14.65"
$ws.Range("E25").Value = "  -6.23%  "
$ws.Range("D26").Value = "3.59"
$ws.Range("E26").Value = "  -3.55%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "2.34"
$ws.Range("E28").Value = "  -8.58%  "
$ws.Range("D29").Value = "31.90"
$ws.Range("E29").Value = "  -6.83%  "
$ws.Range("D30").Value = "8.36"
$ws.Range("E30").Value = "  -7.26%  "
$ws.Range("D31").Value = "7.62"
$ws.Range("E31").Value = "  -6.62%  "
$ws.Range("E32").Value = "  -5.23%  "
$ws.Range("D33").Value = "2.71"
$ws.Range("E33").Value = "  -9.02%  "
$ws.Range("D34").Value = "6.37"
$ws.Range("E34").Value = "  -7.82%  "
$ws.Range("D35").Value = "580.56"
$ws.Range("E35").Value = "  -10.06%  "
$ws.Range("D36").Value = "10.26"
$ws.Range("E36").Value = "  -4.29%  "
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("D38").Value = "0.0932"
$ws.Range("E38").Value = "  -6.49%  "
$ws.Range("D39").Value = "55.84"
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("D40").Value = "0.0455"
$ws.Range("E40").Value = "  -4.27%  "
$ws.Range("D41").Value = "0.139"
$ws.Range("E41").Value = "  -2.55%  "
$ws.Range("D42").Value = "2.95"
$ws.Range("E42").Value = "  -17.88%  "
$ws.Range("D43").Value = "3.176.77"
$ws.Range("E43").Value = "  -5.57%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "30.64"
$ws.Range("E44").Value = "  -5.01%  "
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "0.0₃0652"
$ws.Range("E45").Value = "  -13.10%  "
$ws.Range("D46").Value = "2.73"
$ws.Range("E46").Value = "  -6.40%  "
$ws.Range("D47").Value = "0.286"
$ws.Range("E47").Value = "  -7.92%  "
$ws.Range("D48").Value = "2.31"
$ws.Range("E48").Value = "  -9.51%  "
$ws.Range("D49").Value = "0.124"
$ws.Range("E49").Value = "  -4.55%  "
$ws.Range("D50").Value = "131.50"
$ws.Range("E50").Value = "  -2.14%  "

Write-Host "Applied cryptos update"
